$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural edits: insert a new column C (unit "Name" column) and two new
# rows (4:5) for the "Default From Row" / "Default To Row" config ---
$ws.Columns("C").Insert()
$ws.Rows("4:5").Insert()

# New column C should be the same width as column B (22 chars)
$ws.Columns("C").ColumnWidth = 21.14

# --- New configuration rows 4 and 5 ---
$ws.Range("A4").Value = "Default From Row"
$ws.Range("B4").Value = 8
$ws.Range("A5").Value = "Default To Row"
$ws.Range("B5").Value = 10

# --- New "Name" column values (friendly unit names) for the unit table rows ---
$ws.Range("C8").Value = "meter"
$ws.Range("C9").Value = "millimeter"
$ws.Range("C10").Value = "foot"
$ws.Range("C11").Value = "inch"
$ws.Range("C12").Value = "mile"
$ws.Range("C13").Value = "yard"
$ws.Range("C14").Value = "kilometer"
$ws.Range("C15").Value = "thousandth of an inch"
$ws.Range("C16").Value = "microinch"
$ws.Range("C17").Value = "micron"
$ws.Range("C18").Value = "micrometer"
$ws.Range("C19").Value = "survey mile (US)"
$ws.Range("C20").Value = "survey foot (US)"

# Give the new Name column its own style (Calibri 11, theme color) via a
# dedicated named style, mirroring the workbook's new "Normal 2" cell style.
$normal2 = $wb.Styles.Add("Normal 2")
$normal2.Font.Name = "Calibri"
$normal2.Font.Size = 11

$ws.Range("C8:C20").Style = "Normal 2"
$ws.Range("C8:C20").Font.ThemeColor = 1

# Row 6 (old row 4, "Units" header) and row 7 (old row 5, "Name"/"Factor")
# need their new C cells to pick up the same style as the rest of the row.
$ws.Range("C6").Value = ""
$ws.Range("C7").Value = ""

# Update the selection / active cell to match the saved state
$ws.Range("F7").Select()
